$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tarefas")

# Add new row 58 with the OverFlow/javascript bug entry.
$ws.Cells.Item(58, 1).Value = 42205
$ws.Cells.Item(58, 2).Value = "Concluída"
$ws.Cells.Item(58, 3).Value = "Codificação"
$ws.Cells.Item(58, 4).Value = "OverFlow no javascript, refatoração do Jurassic"

$ws.Range("A58").Select()
